# Add new polling rows (170-174): opinionway and ifop rollings (1/17 and 1/19)
# and cluster17 poll (1/20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=170; Cells=@{
        A=67; B=2022; C=20; D=1; E=16; F="opinionway"; G="online"; H="partially";
        I=1; J=668; K=1; L="T_0.5"; M=10; N=3; O=1; P=5; Q=3; R=24; S=18;
        V=1; W=1; X=18; Y=11; AD=4
    }},
    @{ Row=171; Cells=@{
        A=68; B=2022; C=21; D=1; E=18; F="harris"; G="online"; H="included";
        I=0; J=1833; K=1; L="T_0.5"; M=11; N=2; O=1; P=7; Q=3; R=26; S=16;
        V="T_0.5"; W=1; X=17; Y=14; Z="T_0.5"; AB=1
    }},
    @{ Row=172; Cells=@{
        A=68; B=2022; C=21; D=1; E=18; F="harris"; G="online"; H="included";
        I=0; J=1811; K="T_0.5"; L="T_0.5"; M=10; N=2; O=1; P=6; Q=3; R=25; S=16;
        V="T_0.5"; W=1; X=17; Y=14; Z="T_0.5"; AB=1; AD=4
    }},
    @{ Row=173; Cells=@{
        A=69; B=2022; C=21; D=1; E=17; F="ifop"; G="online"; H="included";
        I=1; J=764; K=1; L=0.5; M=10; N=2; O=0.5; P=5.5; Q=3; R=15; S=16;
        V=1.5; W=1.5; X=18; Y=11.5; AD=3.5; AE=0.5
    }},
    @{ Row=174; Cells=@{
        A=70; B=2022; C=20; D=1; E=13; F="cluster17"; G="online"; H="partially";
        I=0; J=2558; K=1.5; L=0.5; M=12.5; N=2; O=1; P=4.5; Q=2; R=22.5; S=13;
        V=1; W=2.5; X=14.5; Y=14; Z=1.5; AB=1.5; AD=5.5
    }}
)

foreach ($rowSpec in $rows) {
    $r = $rowSpec.Row
    $cells = $rowSpec.Cells
    foreach ($col in $cells.Keys) {
        $ws.Range("$col$r").Value = $cells[$col]
    }
}

# Update the frozen-pane view to match where the author left the cursor.
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A148").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("K173").Select()
